$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Address" reference block added in column K/L (below the existing
# factor tables), matching the pricing sheet's second info block.
$ws.Range("K13").Value = "vehicle year"

$ws.Range("K14").Value = "Type of use "
$ws.Range("L14").Value = "Personal:1 Company:2"

$ws.Range("K15").Value = "Car (neworused)"
$ws.Range("L15").Value = "New:0, used:1"

$ws.Range("K16").Value = "Fuel type"

$ws.Range("K19").Value = "Drivers"
$ws.Range("L19").Value = "Number of drivers, 1,2"

$ws.Range("K20").Value = "Nationality"
$ws.Range("L20").Value = "ID"

# Leave the selection where the author left off editing.
$ws.Range("G16").Select()
